$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Rename the "CART" algorithm row to "DTREE" (row 5, column B)
# ------------------------------------------------------------------
$ws.Range("B5").Value = "DTREE"

# ------------------------------------------------------------------
# 2. Drop the "NB" algorithm row entirely (original row 8). Doing
#    this before the column restructuring keeps the row math simple;
#    the remaining "SVM" row shifts up from row 9 to row 8. Column A
#    holds a literal 0-based index (not a formula), so after the
#    shift it still carries the old literal "7" -- reset it to "6" to
#    match its new position among the seven surviving rows.
# ------------------------------------------------------------------
$ws.Rows("8:8").Delete()
$ws.Range("A8").Value = 6

# ------------------------------------------------------------------
# 3. Expand each of the five "<Horizon> Base" columns (currently
#    B:G, mean values in C:G) into a pair of "<Horizon> Base mean" /
#    "<Horizon> Base std" columns. We insert a new blank column right
#    after each existing metric column -- processing right-to-left so
#    earlier insertions don't disturb the column letters of the
#    horizons still to be processed -- then fill the header + the
#    seven data rows with the refreshed mean and new std values.
# ------------------------------------------------------------------

$horizons = @("One Year Base", "Two Year Base", "Three Year Base", "Five Year Base", "Ten Year Base")
# Original (pre-insert) columns holding the mean values for each horizon.
$meanCols = @("C", "D", "E", "F", "G")
# Column the std column ends up at once it is inserted right after the
# matching mean column above (kept as plain strings -- avoid [char]
# arithmetic, which this host coerces to numeric addition).
$stdCols = @("D", "E", "F", "G", "H")

# Refreshed mean values per horizon, per data row (rows 2-8, i.e. the
# seven surviving algorithms LR, LDA, KNN, DTREE, RTREE, XTREE, SVM).
$meanValues = @{
    "One Year Base"   = @(0.8301764870265979, 0.834042002112146, 0.8188613067032602, 0.7480844923069249, 0.7304882364943679, 0.8226885444408912, 0.8284523723682378)
    "Two Year Base"   = @(0.804474988907784, 0.8101490921849617, 0.7949935364549775, 0.7522325658032371, 0.7170672844357646, 0.8023488690683539, 0.8121991441434286)
    "Three Year Base" = @(0.7899844003189213, 0.784359184675707, 0.7970131973818548, 0.730953756578356, 0.7014888452720773, 0.7876084984814273, 0.7976868946973996)
    "Five Year Base"  = @(0.7602991958696066, 0.758225926686428, 0.7849936330329443, 0.7363863921154541, 0.6857028003138599, 0.7775465523649588, 0.78808773864576)
    "Ten Year Base"   = @(0.7299754826971369, 0.7371086246428523, 0.7501861484506362, 0.7019575220013522, 0.6524342347795187, 0.7409188894969273, 0.76311332558259)
}

# New std values per horizon, per data row.
$stdValues = @{
    "One Year Base"   = @(0.01620696097178251, 0.01864744979041622, 0.01475753081972156, 0.02739293401449815, 0.02373524230449939, 0.01553948102979973, 0.01610605149363529)
    "Two Year Base"   = @(0.01001519155963652, 0.01979147470624407, 0.009704860393678101, 0.02304903023573984, 0.0224952372983046, 0.01688043602760385, 0.01422813002032585)
    "Three Year Base" = @(0.02033703622902071, 0.01907896668677899, 0.01736545885836565, 0.02686600358376494, 0.01459773894989873, 0.02327066474866556, 0.02442913581330806)
    "Five Year Base"  = @(0.03485790028378719, 0.03502610512468584, 0.02155699213976994, 0.03812658847494436, 0.0241335124732715, 0.02852058239499256, 0.02874693741794955)
    "Ten Year Base"   = @(0.02691398258421014, 0.02657984219504058, 0.02194102561853693, 0.03912593830585765, 0.02852221256805043, 0.02674371294456405, 0.03044243359301423)
}

for ($i = $horizons.Length - 1; $i -ge 0; $i--) {
    $horizon = $horizons[$i]
    $meanCol = $meanCols[$i]
    $stdCol = $stdCols[$i]

    # Insert a new column right after the mean column (i.e. before the
    # next column letter); the inserted cells inherit the formatting
    # of their left neighbour, so the new header cell is already bold.
    $ws.Columns($stdCol + ":" + $stdCol).Insert()

    # Header text.
    $ws.Range($meanCol + "1").Value = $horizon + " mean"
    $ws.Range($stdCol + "1").Value = $horizon + " std"

    # Data rows 2-8 (seven surviving algorithms).
    $means = $meanValues[$horizon]
    $stds = $stdValues[$horizon]
    for ($r = 0; $r -lt 7; $r++) {
        $row = $r + 2
        $ws.Range($meanCol + $row).Value = $means[$r]
        $ws.Range($stdCol + $row).Value = $stds[$r]
    }
}

Write-Host "Update complete."
